$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 70.666664
$ws.Range("I11").Value = 70.666664
$ws.Range("K11").Value = 70.666664
$ws.Range("M11").Value = 69.333336
$ws.Range("H17").Value = 2566.8333
$ws.Range("J17").Value = 2566.8333
$ws.Range("L17").Value = 7700.499899999999
$ws.Range("N17").Value = -8036.499899999999
$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3064
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15320
$ws.Range("N77").ClearContents()
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H98").Value = 3463.15
$ws.Range("I98").Value = 3979.1875
$ws.Range("J98").Value = 1399
$ws.Range("K98").Value = 3979.1875
$ws.Range("L98").Value = 1399
$ws.Range("M98").Value = -2481.1875
$ws.Range("N98").Value = -4395
$ws.Range("H122").Value = 3463.15
$ws.Range("I122").Value = 3979.1875
$ws.Range("J122").Value = 1399
$ws.Range("K122").Value = 11937.5625
$ws.Range("L122").Value = 4197
$ws.Range("M122").Value = -9487.5625
$ws.Range("N122").Value = -9097
$ws.Range("H135").Value = 43479004
$ws.Range("I135").Value = 670.0952
$ws.Range("J135").Value = 500001500
$ws.Range("K135").Value = 6030.8568
$ws.Range("L135").Value = 4500013500
$ws.Range("M135").Value = -3495.8568
$ws.Range("N135").Value = -4500018570
$ws.Range("H137").Value = 1282.8478
$ws.Range("I137").Value = 1014.69696
$ws.Range("J137").Value = 1963.5385
$ws.Range("K137").Value = 3044.09088
$ws.Range("L137").Value = 5890.6155
$ws.Range("M137").Value = -494.0908799999997
$ws.Range("N137").Value = -10990.6155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 731.74194
$ws.Range("I2").Value = 431.36
$ws.Range("J2").Value = 1983.3334
$ws.Range("K2").Value = 431.36
$ws.Range("L2").Value = 1983.3334
$ws.Range("M2").Value = -318.36
$ws.Range("N2").Value = -2209.3334
$ws.Range("H32").Value = 4791.45
$ws.Range("I32").Value = 3656.6548
$ws.Range("J32").Value = 10749.125
$ws.Range("K32").Value = 3656.6548
$ws.Range("L32").Value = 10749.125
$ws.Range("M32").Value = -3369.6548
$ws.Range("N32").Value = -11323.125
$ws.Range("H45").Value = 1104
$ws.Range("I45").Value = 1075.7142
$ws.Range("J45").Value = 1170
$ws.Range("K45").Value = 1075.7142
$ws.Range("L45").Value = 1170
$ws.Range("M45").Value = -698.7141999999999
$ws.Range("N45").Value = -1924
$ws.Range("H74").Value = 2392.6843
$ws.Range("I74").Value = 1553.3
$ws.Range("J74").Value = 3325.3333
$ws.Range("K74").Value = 1553.3
$ws.Range("L74").Value = 3325.3333
$ws.Range("M74").Value = -679.3
$ws.Range("N74").Value = -5073.3333
$ws.Range("H77").Value = 2392.6843
$ws.Range("I77").Value = 1553.3
$ws.Range("J77").Value = 3325.3333
$ws.Range("K77").Value = 7766.5
$ws.Range("L77").Value = 16626.6665
$ws.Range("M77").Value = -3398.5
$ws.Range("N77").Value = -25362.6665
$ws.Range("H116").Value = 731.74194
$ws.Range("I116").Value = 431.36
$ws.Range("J116").Value = 1983.3334
$ws.Range("K116").Value = 431.36
$ws.Range("L116").Value = 1983.3334
$ws.Range("M116").Value = 1862.64
$ws.Range("N116").Value = -6571.3334
$ws.Range("H132").Value = 2505.0698
$ws.Range("I132").Value = 2129.7693
$ws.Range("K132").Value = 6389.3079
$ws.Range("M132").Value = -3859.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 731.74194
$ws.Range("I3").Value = 431.36
$ws.Range("J3").Value = 1983.3334
$ws.Range("K3").Value = 431.36
$ws.Range("L3").Value = 1983.3334
$ws.Range("M3").Value = -317.36
$ws.Range("N3").Value = -2211.3334
$ws.Range("H86").Value = 2679.182
$ws.Range("I86").Value = 2638.8096
$ws.Range("J86").Value = 2749.8333
$ws.Range("K86").Value = 2638.8096
$ws.Range("L86").Value = 2749.8333
$ws.Range("M86").Value = -1515.8096
$ws.Range("N86").Value = -4995.8333
$ws.Range("H89").Value = 2679.182
$ws.Range("I89").Value = 2638.8096
$ws.Range("J89").Value = 2749.8333
$ws.Range("K89").Value = 2638.8096
$ws.Range("L89").Value = 2749.8333
$ws.Range("M89").Value = -7578.048000000001
$ws.Range("N89").Value = -24981.1665
$ws.Range("H105").Value = 500001500
$ws.Range("I105").Value = 1000000000
$ws.Range("K105").Value = 1000000000
$ws.Range("M105").Value = -999998253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1318.1428
$ws.Range("I31").Value = 1318.1428
$ws.Range("K31").Value = 1318.1428
$ws.Range("M31").Value = -1023.1428
$ws.Range("H34").Value = 1318.1428
$ws.Range("I34").Value = 1318.1428
$ws.Range("K34").Value = 1318.1428
$ws.Range("M34").Value = -1116.1428
$ws.Range("H134").Value = 16130581
$ws.Range("I134").Value = 1334.2273
$ws.Range("J134").Value = 55557628
$ws.Range("K134").Value = 4002.6819
$ws.Range("L134").Value = 166672884
$ws.Range("M134").Value = -1467.6819
$ws.Range("N134").Value = -166677954

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 4538.077
$ws.Range("J74").Value = 4832.9165
$ws.Range("L74").Value = 14498.7495
$ws.Range("N74").Value = -16620.7495
$ws.Range("H77").Value = 4538.077
$ws.Range("J77").Value = 4832.9165
$ws.Range("L77").Value = 43496.2485
$ws.Range("N77").Value = -54104.2485
$ws.Range("H88").Value = 5373.25
$ws.Range("J88").Value = 5881.3887
$ws.Range("L88").Value = 17644.1661
$ws.Range("N88").Value = -18500.1661
$ws.Range("H91").Value = 5373.25
$ws.Range("J91").Value = 5881.3887
$ws.Range("L91").Value = 17644.1661
$ws.Range("N91").Value = -20608.1661
$ws.Range("H120").Value = 12632.667
$ws.Range("J120").Value = 16449.5
$ws.Range("L120").Value = 49348.5
$ws.Range("N120").Value = -59024.5
$ws.Range("H131").Value = 30306808
$ws.Range("I131").Value = 83333816
$ws.Range("J131").Value = 5661.7144
$ws.Range("K131").Value = 250001448
$ws.Range("L131").Value = 16985.1432
$ws.Range("M131").Value = -249996408
$ws.Range("N131").Value = -27065.1432
$ws.Range("H133").Value = 4141.591
$ws.Range("I133").Value = 2970
$ws.Range("J133").Value = 4486.1763
$ws.Range("K133").Value = 8910
$ws.Range("L133").Value = 13458.5289
$ws.Range("M133").Value = -3850
$ws.Range("N133").Value = -23578.5289
$ws.Range("H137").Value = 22064170
$ws.Range("I137").Value = 57694304
$ws.Range("J137").Value = 7420.5713
$ws.Range("K137").Value = 173082912
$ws.Range("L137").Value = 22261.7139
$ws.Range("M137").Value = -173077812
$ws.Range("N137").Value = -32461.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40912250
$ws.Range("I70").Value = 35717690
$ws.Range("J70").Value = 50002730
$ws.Range("K70").Value = 35717690
$ws.Range("L70").Value = 50002730
$ws.Range("M70").Value = -35717420
$ws.Range("N70").Value = -50003270
$ws.Range("H73").Value = 40912250
$ws.Range("I73").Value = 35717690
$ws.Range("J73").Value = 50002730
$ws.Range("K73").Value = 35717690
$ws.Range("L73").Value = 50002730
$ws.Range("M73").Value = -35716754
$ws.Range("N73").Value = -50004602
$ws.Range("H97").Value = 850
$ws.Range("I97").Value = 700
$ws.Range("K97").Value = 700
$ws.Range("M97").Value = -204
$ws.Range("H102").Value = 1872.091
$ws.Range("I102").Value = 1817
$ws.Range("K102").Value = 1817
$ws.Range("M102").Value = -195
$ws.Range("H113").Value = 1423.1666
$ws.Range("I113").Value = 1277.2142
$ws.Range("J113").Value = 1627.5
$ws.Range("K113").Value = 1277.2142
$ws.Range("L113").Value = 1627.5
$ws.Range("M113").Value = 892.7858000000001
$ws.Range("N113").Value = -5967.5
$ws.Range("H122").Value = 3205.8572
$ws.Range("I122").Value = 3272.0454
$ws.Range("J122").Value = 2963.1667
$ws.Range("K122").Value = 9816.136200000001
$ws.Range("L122").Value = 8889.500100000001
$ws.Range("M122").Value = -7366.136200000001
$ws.Range("N122").Value = -13789.5001
$ws.Range("H126").Value = 1941.3889
$ws.Range("I126").Value = 1694.5555
$ws.Range("J126").Value = 2188.2222
$ws.Range("K126").Value = 5083.666499999999
$ws.Range("L126").Value = 6564.6666
$ws.Range("M126").Value = -2613.666499999999
$ws.Range("N126").Value = -11504.6666
$ws.Range("H132").Value = 3222.8845
$ws.Range("I132").Value = 3032.389
$ws.Range("J132").Value = 3651.5
$ws.Range("K132").Value = 9097.167000000001
$ws.Range("L132").Value = 10954.5
$ws.Range("M132").Value = -6567.167000000001
$ws.Range("N132").Value = -16014.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4521.4614
$ws.Range("I46").Value = 793.3333
$ws.Range("J46").Value = 5639.9
$ws.Range("K46").Value = 793.3333
$ws.Range("L46").Value = 5639.9
$ws.Range("M46").Value = -605.3333
$ws.Range("N46").Value = -6015.9
$ws.Range("H93").Value = 1333.3334
$ws.Range("I93").Value = 1375
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1375
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -127
$ws.Range("N93").Value = -3496
$ws.Range("H136").Value = 2240.5
$ws.Range("I136").Value = 1784.8
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5354.4
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2804.4
$ws.Range("N136").Value = -14100
$ws.Range("H140").Value = 49050.91
$ws.Range("J140").Value = 49050.91
$ws.Range("L140").Value = 49050.91
$ws.Range("N140").Value = -59410.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 14527
$ws.Range("J63").Value = 15427.143
$ws.Range("L63").Value = 15427.143
$ws.Range("N63").Value = -16675.143
$ws.Range("H66").Value = 14527
$ws.Range("J66").Value = 15427.143
$ws.Range("L66").Value = 46281.429
$ws.Range("N66").Value = -52521.429
$ws.Range("H126").Value = 43479144
$ws.Range("I126").Value = 62500496
$ws.Range("J126").Value = 1773.4286
$ws.Range("K126").Value = 187501488
$ws.Range("L126").Value = 5320.2858
$ws.Range("M126").Value = -187499018
$ws.Range("N126").Value = -10260.2858
